# Insert a new weekly price record at row 186 of the "Zapallo italiano" sheet.
# This pushes the existing rows 186-227 down to 187-228 (dimension grows from
# A1:R227 to A1:R228) and fills the newly-opened row 186 with the latest
# observation (Region del Maule, 11-Nov-2021 / serial 44511).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a blank row above the current row 186, shifting rows 186:227 -> 187:228.
$ws.Rows.Item(186).Insert()

$ws.Cells.Item(186, 1).Value = 5
$ws.Cells.Item(186, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(186, 3).Value = 'Maule'
$ws.Cells.Item(186, 4).Value = 44511
$ws.Cells.Item(186, 5).Value = 7
$ws.Cells.Item(186, 6).Value = 100112032
$ws.Cells.Item(186, 7).Value = 'Zapallo italiano'
$ws.Cells.Item(186, 8).Value = 'Sin especificar'
$ws.Cells.Item(186, 9).Value = 'Primera'
$ws.Cells.Item(186, 10).Value = 400
$ws.Cells.Item(186, 11).Value = 7000
$ws.Cells.Item(186, 12).Value = 7000
$ws.Cells.Item(186, 13).Value = 7000
$ws.Cells.Item(186, 14).Value = '$/caja 60 unidades'
$ws.Cells.Item(186, 15).Value = 'Región del Maule'
$ws.Cells.Item(186, 16).Value = 117
$ws.Cells.Item(186, 17).Value = 60
$ws.Cells.Item(186, 18).Value = 'Hortaliza'
